# Clean up Column A "Reference" values: strip the stray trailing "16"
# that was accidentally appended to most verse references (e.g.
# "Mark 1:116" -> "Mark 1:1"). A handful of rows never had the stray
# suffix (e.g. "Mark 13:23", "Mark 16:9") and are left untouched because
# they do not end in "16".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val.EndsWith("16")) {
        $cell.Value = $val.Substring(0, $val.Length - 2)
    }
}
